$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1 and J1 with header style (same as existing header cells like H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data values for I2:J60
$data = @(
    @(2, 8, 8),
    @(3, 7, 7),
    @(4, 5, 7),
    @(5, 9, 9),
    @(6, 5, 6),
    @(7, 9, 9),
    @(8, 1, 2),
    @(9, 7, 8),
    @(10, 8, 9),
    @(11, 9, 10),
    @(12, 7, 7),
    @(13, 8, 8),
    @(14, 1, 1),
    @(15, 1, 2),
    @(16, 1, 1),
    @(17, 12, 13),
    @(18, 6, 6),
    @(19, 1, 2),
    @(20, 7, 8),
    @(21, 7, 8),
    @(22, 1, 1),
    @(23, 1, 2),
    @(24, 10, 10),
    @(25, 1, 1),
    @(26, 1, 2),
    @(27, 1, 3),
    @(28, 9, 9),
    @(29, 7, 8),
    @(30, 1, 1),
    @(31, 9, 9),
    @(32, 8, 8),
    @(33, 7, 8),
    @(34, 8, 8),
    @(35, 1, 1),
    @(36, 8, 8),
    @(37, 6, 8),
    @(38, 3, 5),
    @(39, 7, 8),
    @(40, 7, 7),
    @(41, 6, 7),
    @(42, 6, 7),
    @(43, 6, 8),
    @(44, 5, 7),
    @(45, 7, 8),
    @(46, 5, 5),
    @(47, 5, 6),
    @(48, 7, 9),
    @(49, 9, 9),
    @(50, 6, 8),
    @(51, 3, 5),
    @(52, 8, 9),
    @(53, 9, 9),
    @(54, 8, 9),
    @(55, 5, 6),
    @(56, 5, 6),
    @(57, 7, 8),
    @(58, 3, 4),
    @(59, 6, 7),
    @(60, 4, 5)
)

foreach ($row in $data) {
    $r = $row[0]
    $iVal = $row[1]
    $jVal = $row[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}

Write-Output "Added I0/IF columns for rows 1-60"